$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.25418433921579
$ws.Range("C2").Value = 16.09581905353889
$ws.Range("D2").Value = 14.35991293450656
$ws.Range("E2").Value = 14.93513336620661
$ws.Range("G2").Value = 3.793155217024867
$ws.Range("I2").Value = 37.6659530161923
$ws.Range("J2").Value = 8.82351175058921
$ws.Range("K2").Value = 19.12735587125486
$ws.Range("M2").Value = 21.98481685026091
$ws.Range("N2").Value = 26.43419441940317
$ws.Range("B3").Value = 17.12391956109014
$ws.Range("C3").Value = 15.98405259244906
$ws.Range("D3").Value = 14.3444761988353
$ws.Range("E3").Value = 14.94076598569405
$ws.Range("G3").Value = 3.796767125185716
$ws.Range("I3").Value = 37.61144220260938
$ws.Range("J3").Value = 8.840621077749718
$ws.Range("K3").Value = 19.04235823541827
$ws.Range("M3").Value = 21.9647693011153
$ws.Range("N3").Value = 26.44213538842156
$ws.Range("B4").Value = 17.04841928110624
$ws.Range("C4").Value = 15.91919679858167
$ws.Range("D4").Value = 14.33790738031995
$ws.Range("E4").Value = 14.94685528084141
$ws.Range("G4").Value = 3.799099635494282
$ws.Range("I4").Value = 37.58320638954465
$ws.Range("J4").Value = 8.851896124568841
$ws.Range("K4").Value = 18.99481787523009
$ws.Range("M4").Value = 21.95699638179263
$ws.Range("N4").Value = 26.44872963585438
$ws.Range("B5").Value = 17.01880893173174
$ws.Range("C5").Value = 15.8937367004669
$ws.Range("D5").Value = 14.33596366010755
$ws.Range("E5").Value = 14.94999822028002
$ws.Range("G5").Value = 3.800079124422488
$ws.Range("I5").Value = 37.57301963624448
$ws.Range("J5").Value = 8.856684772901444
$ws.Range("K5").Value = 18.97662764319254
$ws.Range("M5").Value = 21.954971751991
$ws.Range("N5").Value = 26.45184879563389
$ws.Range("B6").Value = 17.01396282790125
$ws.Range("C6").Value = 15.88956819467926
$ws.Range("D6").Value = 14.33568522124235
$ws.Range("E6").Value = 14.95056005037477
$ws.Range("G6").Value = 3.800243520843365
$ws.Range("I6").Value = 37.57140791723484
$ws.Range("J6").Value = 8.857491651423189
$ws.Range("K6").Value = 18.97367898518155
$ws.Range("M6").Value = 21.95470464178401
$ws.Range("N6").Value = 26.45239281272828
$ws.Range("B7").Value = 17.04801522533345
$ws.Range("C7").Value = 15.91884948382841
$ws.Range("D7").Value = 14.3378781965598
$ws.Range("E7").Value = 14.94689498945422
$ws.Range("G7").Value = 3.799112727775211
$ws.Range("I7").Value = 37.58306366074464
$ws.Range("J7").Value = 8.851959920043145
$ws.Range("K7").Value = 18.99456774894833
$ws.Range("M7").Value = 21.95696444706174
$ws.Range("N7").Value = 26.44876995332576
$ws.Range("B8").Value = 17.20835693624498
$ws.Range("C8").Value = 16.05651300817981
$ws.Range("D8").Value = 14.35398756193255
$ws.Range("E8").Value = 14.93652949632118
$ws.Range("G8").Value = 3.794376845051008
$ws.Range("I8").Value = 37.64607362899816
$ws.Range("J8").Value = 8.829251498197568
$ws.Range("K8").Value = 19.09709273952664
$ws.Range("M8").Value = 21.97696446023685
$ws.Range("N8").Value = 26.43657570249807
$ws.Range("B9").Value = 17.5569614656122
$ws.Range("C9").Value = 16.35535932541371
$ws.Range("D9").Value = 14.40857636389351
$ws.Range("E9").Value = 14.93707228650118
$ws.Range("G9").Value = 3.78599559732658
$ws.Range("I9").Value = 37.81098952931313
$ws.Range("J9").Value = 8.79081143326342
$ws.Range("K9").Value = 19.33434103569485
$ws.Range("M9").Value = 22.05205344251895
$ws.Range("N9").Value = 26.42630669269455
$ws.Range("B10").Value = 17.8319166593815
$ws.Range("C10").Value = 16.59106350668357
$ws.Range("D10").Value = 14.46255974154352
$ws.Range("E10").Value = 14.95017634467596
$ws.Range("G10").Value = 3.780383159583796
$ws.Range("I10").Value = 37.95708991995534
$ws.Range("J10").Value = 8.76625918744662
$ws.Range("K10").Value = 19.52969169380732
$ws.Range("M10").Value = 22.12887888691689
$ws.Range("N10").Value = 26.42709312558734
$ws.Range("B11").Value = 17.96063321290007
$ws.Range("C11").Value = 16.70145628853206
$ws.Range("D11").Value = 14.49009075897861
$ws.Range("E11").Value = 14.95888805863527
$ws.Range("G11").Value = 3.777946833846719
$ws.Range("I11").Value = 38.02889558395506
$ws.Range("J11").Value = 8.75588595165909
$ws.Range("K11").Value = 19.62288380068163
$ws.Range("M11").Value = 22.1684716600693
$ws.Range("N11").Value = 26.42926221462374
$ws.Range("B12").Value = 18.00985784194643
$ws.Range("C12").Value = 16.74368478214142
$ws.Range("D12").Value = 14.50093935236318
$ws.Range("E12").Value = 14.96258128875459
$ws.Range("G12").Value = 3.777040942714673
$ws.Range("I12").Value = 38.05684718586534
$ws.Range("J12").Value = 8.75207191605058
$ws.Range("K12").Value = 19.65877247472243
$ws.Range("M12").Value = 22.18412591894965
$ws.Range("N12").Value = 26.43034410937786
$ws.Range("B13").Value = 17.99923563562171
$ws.Range("C13").Value = 16.73457170864345
$ws.Range("D13").Value = 14.49858416752915
$ws.Range("E13").Value = 14.96176837176832
$ws.Range("G13").Value = 3.777235301810991
$ws.Range("I13").Value = 38.0507936281756
$ws.Range("J13").Value = 8.75288826829747
$ws.Range("K13").Value = 19.65101695476039
$ws.Range("M13").Value = 22.18072519307772
$ws.Range("N13").Value = 26.43009951764313
$ws.Range("B14").Value = 17.96467351927625
$ws.Range("C14").Value = 16.70492211004521
$ws.Range("D14").Value = 14.49097482488102
$ws.Range("E14").Value = 14.95918400869909
$ws.Range("G14").Value = 3.777871971683301
$ws.Range("I14").Value = 38.03117998947742
$ws.Range("J14").Value = 8.755569884169443
$ws.Range("K14").Value = 19.62582451885202
$ws.Range("M14").Value = 22.16974633129731
$ws.Range("N14").Value = 26.42934600172607
$ws.Range("B15").Value = 17.94356483508795
$ws.Range("C15").Value = 16.68681531734827
$ws.Range("D15").Value = 14.48636885572858
$ws.Range("E15").Value = 14.9576523190858
$ws.Range("G15").Value = 3.778264121171454
$ws.Range("I15").Value = 38.01926483697585
$ws.Range("J15").Value = 8.757227298419675
$ws.Range("K15").Value = 19.61047070277779
$ws.Range("M15").Value = 22.1631073801656
$ws.Range("N15").Value = 26.42891837768985
$ws.Range("B16").Value = 17.82357447806837
$ws.Range("C16").Value = 16.58391025989121
$ws.Range("D16").Value = 14.4608199966245
$ws.Range("E16").Value = 14.94966226274058
$ws.Range("G16").Value = 3.780544720689734
$ws.Range("I16").Value = 37.95250410964527
$ws.Range("J16").Value = 8.766953085459708
$ws.Range("K16").Value = 19.52368646204916
$ws.Range("M16").Value = 22.12638432321546
$ws.Range("N16").Value = 26.42698782364787
$ws.Range("B17").Value = 17.75086740134801
$ws.Range("C17").Value = 16.52157164223731
$ws.Range("D17").Value = 14.44590524667809
$ws.Range("E17").Value = 14.9454643175648
$ws.Range("G17").Value = 3.781973636421066
$ws.Range("I17").Value = 37.91291175119692
$ws.Range("J17").Value = 8.77312309833377
$ws.Range("K17").Value = 19.47153899227187
$ws.Range("M17").Value = 22.10504106755569
$ws.Range("N17").Value = 26.42626746029343
$ws.Range("B18").Value = 17.7093931421499
$ws.Range("C18").Value = 16.48601636574113
$ws.Range("D18").Value = 14.43760687243769
$ws.Range("E18").Value = 14.94330876176155
$ws.Range("G18").Value = 3.782806510900843
$ws.Range("I18").Value = 37.89064279733472
$ws.Range("J18").Value = 8.776746841959476
$ws.Range("K18").Value = 19.44195380929346
$ws.Range("M18").Value = 22.09320259196774
$ws.Range("N18").Value = 26.42602361841552
$ws.Range("B19").Value = 17.69541113806511
$ws.Range("C19").Value = 16.47403044858459
$ws.Range("D19").Value = 14.43484543603507
$ws.Range("E19").Value = 14.94262344321481
$ws.Range("G19").Value = 3.783090400290615
$ws.Range("I19").Value = 37.88318962639708
$ws.Range("J19").Value = 8.777986656415218
$ws.Range("K19").Value = 19.43200765215413
$ws.Range("M19").Value = 22.08926963286447
$ws.Range("N19").Value = 26.42597033794185
$ws.Range("B20").Value = 17.75857178723124
$ws.Range("C20").Value = 16.52817683964282
$ws.Range("D20").Value = 14.44746397860625
$ws.Range("E20").Value = 14.94588439959708
$ws.Range("G20").Value = 3.781820388202583
$ws.Range("I20").Value = 37.91707436895121
$ws.Range("J20").Value = 8.77245853894569
$ws.Range("K20").Value = 19.47704804931733
$ws.Range("M20").Value = 22.10726784866774
$ws.Range("N20").Value = 26.42632649864348
$ws.Range("B21").Value = 17.97481248481762
$ws.Range("C21").Value = 16.71361962587626
$ws.Range("D21").Value = 14.49319842552771
$ws.Range("E21").Value = 14.95993240893992
$ws.Range("G21").Value = 3.777684514098165
$ws.Range("I21").Value = 38.0369204230344
$ws.Range("J21").Value = 8.754779135060042
$ws.Range("K21").Value = 19.63320808623858
$ws.Range("M21").Value = 22.17295319644421
$ws.Range("N21").Value = 26.42956025740623
$ws.Range("B22").Value = 18.11893108453795
$ws.Range("C22").Value = 16.83728053122012
$ws.Range("D22").Value = 14.52555284222065
$ws.Range("E22").Value = 14.97141108070508
$ws.Range("G22").Value = 3.775078730667698
$ws.Range("I22").Value = 38.11967438934569
$ws.Range("J22").Value = 8.743889435261396
$ws.Range("K22").Value = 19.7387466853107
$ws.Range("M22").Value = 22.21973358606067
$ws.Range("N22").Value = 26.4331921318154
$ws.Range("B23").Value = 18.04177036324376
$ws.Range("C23").Value = 16.77106530779
$ws.Range("D23").Value = 14.50806078155038
$ws.Range("E23").Value = 14.96507495368775
$ws.Range("G23").Value = 3.776460621979018
$ws.Range("I23").Value = 38.07510478004999
$ws.Range("J23").Value = 8.74964075321882
$ws.Range("K23").Value = 19.68210838284653
$ws.Range("M23").Value = 22.19441597418911
$ws.Range("N23").Value = 26.43111479423322
$ws.Range("B24").Value = 17.75508761580165
$ws.Range("C24").Value = 16.52518974287763
$ws.Range("D24").Value = 14.44675841490651
$ws.Range("E24").Value = 14.94569367710879
$ws.Range("G24").Value = 3.781889636290789
$ws.Range("I24").Value = 37.91519091244913
$ws.Range("J24").Value = 8.772758747864266
$ws.Range("K24").Value = 19.47455617219926
$ws.Range("M24").Value = 22.10625977466287
$ws.Range("N24").Value = 26.42629927686907
$ws.Range("B25").Value = 17.45919577536643
$ws.Range("C25").Value = 16.27156692013156
$ws.Range("D25").Value = 14.39135933674146
$ws.Range("E25").Value = 14.93469204556397
$ws.Range("G25").Value = 3.788166693239797
$ws.Range("I25").Value = 37.76197330375021
$ws.Range("J25").Value = 8.800560851360778
$ws.Range("K25").Value = 19.2663836184476
$ws.Range("M25").Value = 22.02791845344457
$ws.Range("N25").Value = 26.42762251324469
